# Update the "取得日時" (retrieved datetime) timestamps in column A of the
# "ランサーズ" sheet for the newly appended batch of rows (rows 2-8) from
# 2025-12-08 01:22:52 to 2025-12-08 01:54:31.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-12-08 01:22:52"
$newTimestamp = "2025-12-08 01:54:31"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
